# Apply updated values from R script results to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K8").Value = 3731.79

$ws.Range("N9").Value = 2187.18
$ws.Range("O9").Value = 2187.18

$ws.Range("M10").Value = 14028.84
$ws.Range("N10").Value = 2680.33
$ws.Range("O10").Value = 2680.33

$ws.Range("N12").Value = 370774.37
$ws.Range("O12").Value = 366033.83

$ws.Range("M13").Value = 208688.45
$ws.Range("N13").Value = 122572.39
$ws.Range("O13").Value = 122541.11

$ws.Range("K14").Value = 0
$ws.Range("N14").Value = 20872.28
$ws.Range("O14").Value = 20872.28

$ws.Range("M15").Value = 65.13
$ws.Range("N15").Value = 61.65
$ws.Range("O15").Value = 61.65

$ws.Range("K19").Value = 14553.33

$ws.Range("K20").Value = 7040.23
$ws.Range("N20").Value = 4784.2
$ws.Range("O20").Value = 4784.2

$ws.Range("K22").Value = 0

$ws.Range("N23").Value = 567.67
$ws.Range("O23").Value = 567.67

$ws.Range("K24").Value = 165086.12

$ws.Range("N26").Value = 46580
$ws.Range("O26").Value = 46080
